$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds "no"/"yes" execute flags for rows 2-11 (row 1 is the header).
# Rows 2-10 currently read "no"; flip them to "yes" (row 11 is already "yes").
# Once no cell references the "no" shared string any more it drops out of
# the shared-strings table automatically (uniqueCount 25 -> 24).
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = "yes"
}

# Scroll the sheet view so column B becomes the left-most visible column.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2

